$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide with the "project hosted at" bullet: update the lfs-course repo
# address from the old "openeuler-practice-courses" org to "openeuler".
# ---------------------------------------------------------------------------
$oldUrl = "https://gitee.com/openeuler-practice-courses/lfs-course"
$newUrl = "https://gitee.com/openeuler/lfs-course"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $shp = $s.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            $full = $tr.Text
            if ($full.Contains($oldUrl)) {
                $idx = $full.IndexOf($oldUrl)
                while ($idx -ge 0) {
                    $sub = $tr.Characters($idx + 1, $oldUrl.Length)
                    $sub.Text = $newUrl
                    $full = $tr.Text
                    $idx = $full.IndexOf($oldUrl)
                }
            }
        }
    }
}
